$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.817.29"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "2.356.83"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.00"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.669"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.28"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.79"
$ws.Range("E11").Value = "  +4.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "37.05"
$ws.Range("E12").Value = "  +14.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.32"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "2.710.93"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.39"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.924"
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").Value = "2.363.00"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "43.777.19"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.60"
$ws.Range("E21").Value = "  -4.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.24"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.51"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  +3.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.86"
$ws.Range("E26").Value = "  -4.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.50"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.64"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.32"
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.98"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.129"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.133"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0754"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.13"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.80"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.60"
$ws.Range("E38").Value = "  +4.42%  "
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0279"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.59"
$ws.Range("E41").Value = "  +17.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.79"
$ws.Range("E42").Value = "  +10.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.77"
$ws.Range("E43").Value = "  +10.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.107"
$ws.Range("E44").Value = "  -3.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.06"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.201"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.21"
$ws.Range("E51").Value = "  -2.54%  "
